# "Generate Report for Handoff"
# The localization status report moves from "In Translation" to
# "Ready for handoff": update the Status text (Overview sheet's per-language
# status columns + each language sheet's Status column) and refresh the
# associated handoff timestamps. Excel widens the Status columns to fit the
# longer text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$ws1.Range("E2").Value = "Ready for handoff"   # Overview: zh-cn status
$ws1.Range("F2").Value = "Ready for handoff"   # Overview: de-de status
$ws2.Range("C2").Value = "Ready for handoff"   # zh-cn sheet: Status
$ws3.Range("C2").Value = "Ready for handoff"   # de-de sheet: Status

# --- Timestamps bumped as part of the handoff refresh ---------------------
# Overview "Latest HO Xliff Generate Date" / de-de "Latest Handback DateTime"
$ws1.Range("G2").Value = "2016-09-01 16:47:09"
$ws3.Range("H2").Value = "2016-09-01 16:47:09"

# zh-cn "Latest Handoff Datetime"
$ws2.Range("H2").Value = "2016-09-01 16:46:58"

# --- Column widths: Excel auto-widens the Status columns for the longer
# text ("Ready for handoff" vs "In Translation"). ---------------------------
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332   # Overview col E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332   # Overview col F (de-de status)
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332   # zh-cn sheet col C (Status)
$ws3.Columns.Item(3).ColumnWidth = 16.333333333333332   # de-de sheet col C (Status)
